# Auto-generated Excel COM-interop script applying cell updates
# described by the upstream diff to Sheets/Atomos_Profits.xlsx
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 1244.7059
$ws.Range("I33").Value = 984
$ws.Range("J33").Value = 1538
$ws.Range("K33").Value = 984
$ws.Range("L33").Value = 1538
$ws.Range("M33").Value = -755
$ws.Range("N33").Value = -1996

$ws.Range("H58").Value = 92285
$ws.Range("I58").Value = 265
$ws.Range("K58").Value = 795
$ws.Range("M58").Value = -645

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 29123.2
$ws.Range("I12").Value = 600
$ws.Range("J12").Value = 36254
$ws.Range("K12").Value = 600
$ws.Range("L12").Value = 36254
$ws.Range("M12").Value = -427
$ws.Range("N12").Value = -36600

$ws.Range("H74").Value = 1519.5135
$ws.Range("I74").Value = 1027.3
$ws.Range("J74").Value = 3629
$ws.Range("K74").Value = 1027.3
$ws.Range("L74").Value = 3629
$ws.Range("M74").Value = -153.3
$ws.Range("N74").Value = -5377

$ws.Range("H77").Value = 1519.5135
$ws.Range("I77").Value = 1027.3
$ws.Range("J77").Value = 3629
$ws.Range("K77").Value = 5136.5
$ws.Range("L77").Value = 18145
$ws.Range("M77").Value = -768.5
$ws.Range("N77").Value = -26881

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 20235.5
$ws.Range("I26").Value = 20235.5
$ws.Range("K26").Value = 20235.5
$ws.Range("M26").Value = -19943.5

$ws.Range("H58").Value = 47000
$ws.Range("J58").Value = 47000
$ws.Range("L58").Value = 47000
$ws.Range("N58").Value = -47588

$ws.Range("H59").Value = 29666.666
$ws.Range("J59").Value = 29666.666
$ws.Range("L59").Value = 29666.666
$ws.Range("N59").Value = -31360.666

$ws.Range("H99").Value = 4219.8945
$ws.Range("I99").Value = 3463.0833
$ws.Range("J99").Value = 5517.2856
$ws.Range("K99").Value = 3463.0833
$ws.Range("L99").Value = 5517.2856
$ws.Range("M99").Value = -1965.0833
$ws.Range("N99").Value = -8513.285599999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 51793.332
$ws.Range("I10").Value = 1525
$ws.Range("J10").Value = 92008
$ws.Range("K10").Value = 1525
$ws.Range("L10").Value = 92008
$ws.Range("M10").Value = -1386
$ws.Range("N10").Value = -92286

$ws.Range("H31").Value = 4172225
$ws.Range("I31").Value = 8337218.5
$ws.Range("J31").Value = 7231.1665
$ws.Range("K31").Value = 8337218.5
$ws.Range("L31").Value = 7231.1665
$ws.Range("M31").Value = -8336923.5
$ws.Range("N31").Value = -7821.1665

$ws.Range("H34").Value = 4172225
$ws.Range("I34").Value = 8337218.5
$ws.Range("J34").Value = 7231.1665
$ws.Range("K34").Value = 8337218.5
$ws.Range("L34").Value = 7231.1665
$ws.Range("M34").Value = -8337016.5
$ws.Range("N34").Value = -7635.1665

$ws.Range("H58").Value = 9618726
$ws.Range("I58").Value = 1633.931
$ws.Range("J58").Value = 21744626
$ws.Range("K58").Value = 1633.931
$ws.Range("L58").Value = 21744626
$ws.Range("M58").Value = -1430.931
$ws.Range("N58").Value = -21745032

$ws.Range("H105").Value = 2826.7917
$ws.Range("I105").Value = 2743.9412
$ws.Range("J105").Value = 3028
$ws.Range("K105").Value = 2743.9412
$ws.Range("L105").Value = 3028
$ws.Range("M105").Value = -996.9412000000002
$ws.Range("N105").Value = -6522

$ws.Range("H131").Value = 26590.285
$ws.Range("I131").Value = 10300
$ws.Range("J131").Value = 27843.385
$ws.Range("K131").Value = 10300
$ws.Range("L131").Value = 27843.385
$ws.Range("M131").Value = -5260
$ws.Range("N131").Value = -37923.38499999999

$ws.Range("H136").Value = 9618726
$ws.Range("I136").Value = 1633.931
$ws.Range("J136").Value = 21744626
$ws.Range("K136").Value = 4901.793
$ws.Range("L136").Value = 65233878
$ws.Range("M136").Value = -2351.793
$ws.Range("N136").Value = -65238978

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 282.85715
$ws.Range("J7").Value = 195
$ws.Range("L7").Value = 585
$ws.Range("N7").Value = -809

$ws.Range("H63").Value = 8186.6665
$ws.Range("I63").Value = 2045
$ws.Range("J63").Value = 13100
$ws.Range("K63").Value = 6135
$ws.Range("L63").Value = 39300
$ws.Range("M63").Value = -5386
$ws.Range("N63").Value = -40798

$ws.Range("H64").Value = 12402.667
$ws.Range("I64").Value = 562
$ws.Range("J64").Value = 27203.5
$ws.Range("K64").Value = 1686
$ws.Range("L64").Value = 81610.5
$ws.Range("M64").Value = -1416
$ws.Range("N64").Value = -82150.5

$ws.Range("H66").Value = 8186.6665
$ws.Range("I66").Value = 2045
$ws.Range("J66").Value = 13100
$ws.Range("K66").Value = 18405
$ws.Range("L66").Value = 117900
$ws.Range("M66").Value = -14661
$ws.Range("N66").Value = -125388

$ws.Range("H67").Value = 12402.667
$ws.Range("I67").Value = 562
$ws.Range("J67").Value = 27203.5
$ws.Range("K67").Value = 1686
$ws.Range("L67").Value = 81610.5
$ws.Range("M67").Value = -750
$ws.Range("N67").Value = -83482.5

$ws.Range("H82").Value = 3600
$ws.Range("I82").Value = 2000
$ws.Range("K82").Value = 6000
$ws.Range("M82").Value = -5594

$ws.Range("H85").Value = 3600
$ws.Range("I85").Value = 2000
$ws.Range("K85").Value = 6000
$ws.Range("M85").Value = -4596

$ws.Range("H88").Value = 4000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 4000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 12000
$ws.Range("M88").Value = ""
$ws.Range("N88").Value = -12856

$ws.Range("H91").Value = 4000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 4000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 12000
$ws.Range("M91").Value = ""
$ws.Range("N91").Value = -14964

$ws.Range("H113").Value = 1695638.5
$ws.Range("I113").Value = 14286239
$ws.Range("J113").Value = 749.9808
$ws.Range("K113").Value = 42858717
$ws.Range("L113").Value = 2249.9424
$ws.Range("M113").Value = -42856547
$ws.Range("N113").Value = -6589.9424

$ws.Range("H123").Value = 2980
$ws.Range("I123").Value = 1450
$ws.Range("K123").Value = 4350
$ws.Range("M123").Value = -1900

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 49221.684
$ws.Range("I102").Value = 3358.1177
$ws.Range("K102").Value = 3358.1177
$ws.Range("M102").Value = -1736.1177

$ws.Range("H107").Value = 994.4
$ws.Range("I107").Value = 411.55554
$ws.Range("J107").Value = 1471.2727
$ws.Range("K107").Value = 411.55554
$ws.Range("L107").Value = 1471.2727
$ws.Range("M107").Value = 1508.44446
$ws.Range("N107").Value = -5311.2727

$ws.Range("H113").Value = 2190.2144
$ws.Range("I113").Value = 1995.5555
$ws.Range("J113").Value = 2540.6
$ws.Range("K113").Value = 1995.5555
$ws.Range("L113").Value = 2540.6
$ws.Range("M113").Value = 174.4445000000001
$ws.Range("N113").Value = -6880.6

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2776
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 2776
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 2776
$ws.Range("M22").Value = ""
$ws.Range("N22").Value = -3366

$ws.Range("H24").Value = 100007
$ws.Range("J24").Value = 100007
$ws.Range("L24").Value = 100007
$ws.Range("N24").Value = -100693

$ws.Range("H27").Value = 2776
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 2776
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 2776
$ws.Range("M27").Value = ""
$ws.Range("N27").Value = -2990

$ws.Range("H46").Value = 475
$ws.Range("I46").Value = 475
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 475
$ws.Range("L46").Value = 0
$ws.Range("N46").Value = ""
$ws.Range("M46").Value = -287

$ws.Range("H122").Value = 3678.5557
$ws.Range("I122").Value = 2701.2727
$ws.Range("K122").Value = 8103.8181
$ws.Range("M122").Value = -5653.8181

$ws.Range("H136").Value = 2178126
$ws.Range("I136").Value = 3228799.5
$ws.Range("J136").Value = 6734
$ws.Range("K136").Value = 9686398.5
$ws.Range("L136").Value = 20202
$ws.Range("M136").Value = -9683848.5
$ws.Range("N136").Value = -25302

$ws.Range("H137").Value = 29000
$ws.Range("J137").Value = 29000
$ws.Range("L137").Value = 29000
$ws.Range("N137").Value = -39200

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 10692.333
$ws.Range("J41").Value = 10692.333
$ws.Range("L41").Value = 10692.333
$ws.Range("N41").Value = -11472.333
